$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.032.36"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "3.588.13"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.53"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.01"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -1.13%  "
$ws.Range("D8").Value = "3.582.33"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +4.18%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.89"
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("E13").Value = "  +7.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.69"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "4.167.68"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.02"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").Value = "3.591.39"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "70.024.67"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "480.47"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.13"
$ws.Range("E23").Value = "  +9.84%  "
$ws.Range("E24").Value = "  -6.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.40"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.95"
$ws.Range("E26").Value = "  +6.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.01"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.16"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.65"
$ws.Range("E31").Value = "  +2.70%  "
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.74"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "583.86"
$ws.Range("E35").Value = "  -6.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.01"
$ws.Range("E36").Value = "  +3.42%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.398"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.27"
$ws.Range("E40").Value = "  +22.77%  "
$ws.Range("E41").Value = "  -3.68%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.136"
$ws.Range("E42").Value = "  -6.46%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.218.14"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("E44").Value = "  +7.14%  "
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0450"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("E47").Value = "  +5.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.13"
$ws.Range("E51").Value = "  -5.37%  "
